# Update "想去人数" (interest count) figures in the F column for the
# "展览" and "全部类型" sheets, reflecting a fresh scrape of source data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 791
$ws1.Range("F5").Value = 150
$ws1.Range("F6").Value = 19
$ws1.Range("F7").Value = 175
$ws1.Range("F8").Value = 358
$ws1.Range("F9").Value = 470
$ws1.Range("F11").Value = 148
$ws1.Range("F12").Value = 11996
$ws1.Range("F13").Value = 5438

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 791
$ws4.Range("F7").Value = 150
$ws4.Range("F8").Value = 19
$ws4.Range("F9").Value = 175
$ws4.Range("F10").Value = 358
$ws4.Range("F11").Value = 470
$ws4.Range("F13").Value = 148
$ws4.Range("F14").Value = 11996
$ws4.Range("F16").Value = 5438
